$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so values like "30.379.53" or
# "1.001" are preserved exactly as literal text (matching the source data)
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.379.53"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "1.872.52"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "237.93"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "0.4816"
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("D8").Value = "0.2815"
$ws.Range("E8").Value = "  -2.14%  "

$ws.Range("D9").Value = "0.06521"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").Value = "1.877.56"
$ws.Range("E10").Value = "  -0.12%  "

$ws.Range("D11").Value = "0.07463"
$ws.Range("E11").Value = "  +2.35%  "

$ws.Range("D12").Value = "16.48"
$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("D13").Value = "5.080"
$ws.Range("E13").Value = "  -2.65%  "

$ws.Range("D14").Value = "87.94"
$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").Value = "0.6560"
$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("D16").Value = "30.358.33"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "13.30"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").Value = "0.000007616"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").Value = "2.122.70"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "222.47"
$ws.Range("E22").Value = "  +14.23%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.283"
$ws.Range("E23").Value = "  -0.31%  "

$ws.Range("D24").Value = "6.184"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").Value = "9.241"
$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").Value = "165.81"
$ws.Range("E26").Value = "  +4.00%  "

$ws.Range("D27").Value = "18.59"
$ws.Range("E27").Value = "  +3.07%  "

$ws.Range("D28").Value = "1.975"
$ws.Range("E28").Value = "  +3.37%  "

$ws.Range("D29").Value = "1.456"
$ws.Range("E29").Value = "  +1.16%  "

$ws.Range("D30").Value = "0.09411"
$ws.Range("E30").Value = "  +2.98%  "

$ws.Range("D31").Value = "4.311"
$ws.Range("E31").Value = "  +1.18%  "

$ws.Range("D32").Value = "4.008"
$ws.Range("E32").Value = "  -1.35%  "

$ws.Range("D33").Value = "0.05036"
$ws.Range("E33").Value = "  -1.47%  "

$ws.Range("D34").Value = "1.212"
$ws.Range("E34").Value = "  +10.75%  "

$ws.Range("D35").Value = "0.7496"
$ws.Range("E35").Value = "  +4.67%  "

$ws.Range("D36").Value = "2.700"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").Value = "0.01836"
$ws.Range("E37").Value = "  +2.29%  "

$ws.Range("D38").Value = "2.618"
$ws.Range("E38").Value = "  -0.76%  "

$ws.Range("D39").Value = "2.078"
$ws.Range("E39").Value = "  +2.03%  "

$ws.Range("D40").Value = "0.9060"
$ws.Range("E40").Value = "  -1.42%  "

$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").Value = "106.91"
$ws.Range("E41").Value = "  +0.59%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "5.917"
$ws.Range("E42").Value = "  +2.08%  "

$ws.Range("D43").Value = "0.4283"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("D45").Value = "7.438"
$ws.Range("E45").Value = "  +0.70%  "

$ws.Range("D46").Value = "64.99"
$ws.Range("E46").Value = "  -1.81%  "

$ws.Range("D47").Value = "0.1299"
$ws.Range("E47").Value = "  -1.61%  "

$ws.Range("D48").Value = "1.493"
$ws.Range("E48").Value = "  +9.21%  "

$ws.Range("D49").Value = "8.996"
$ws.Range("E49").Value = "  -0.77%  "

$ws.Range("D50").Value = "34.22"
$ws.Range("E50").Value = "  +0.79%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05662"
$ws.Range("E51").Value = "  -1.60%  "
